$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.604.22"
$ws.Range("E2").Value = "  +3.26%  "
$ws.Range("D3").Value = "1.696.52"
$ws.Range("E3").Value = "  +1.99%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'316.91"
$ws.Range("E5").Value = "  +2.43%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "'0.3949"
$ws.Range("E7").Value = "  +1.83%  "
$ws.Range("D8").Value = "'0.4016"
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("D9").Value = "'1.526"
$ws.Range("E9").Value = "  +5.02%  "
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").Value = "'52.45"
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("D12").Value = "'0.08769"
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("D13").Value = "'7.242"
$ws.Range("E13").Value = "  +6.78%  "
$ws.Range("D14").Value = "'23.28"
$ws.Range("E14").Value = "  +2.75%  "
$ws.Range("D15").Value = "'8.214"
$ws.Range("E15").Value = "  +12.14%  "
$ws.Range("D16").Value = "'0.00001314"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").Value = "1.695.17"
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("D18").Value = "'99.92"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").Value = "'0.07076"
$ws.Range("E19").Value = "  +2.75%  "
$ws.Range("D20").Value = "'19.68"
$ws.Range("E20").Value = "  +3.19%  "
$ws.Range("D21").Value = "'7.025"
$ws.Range("E21").Value = "  +5.97%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'14.19"
$ws.Range("E23").Value = "  +2.84%  "
$ws.Range("D24").Value = "24.607.23"
$ws.Range("E24").Value = "  +3.30%  "
$ws.Range("D25").Value = "'3.114"
$ws.Range("E25").Value = "  +9.07%  "
$ws.Range("D26").Value = "'2.338"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("D27").Value = "'22.85"
$ws.Range("E27").Value = "  +5.16%  "
$ws.Range("D28").Value = "'162.11"
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("D29").Value = "'136.47"
$ws.Range("E29").Value = "  +4.68%  "
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("D31").Value = "'7.504"
$ws.Range("E31").Value = "  +9.65%  "
$ws.Range("D32").Value = "1.882.40"
$ws.Range("D33").Value = "'1.082"
$ws.Range("E33").Value = "  -2.93%  "
$ws.Range("D34").Value = "'0.08585"
$ws.Range("E34").Value = "  +0.46%  "
$ws.Range("D35").Value = "'7.179"
$ws.Range("E35").Value = "  +6.67%  "
$ws.Range("D36").Value = "'11.53"
$ws.Range("E36").Value = "  +10.37%  "
$ws.Range("D37").Value = "'0.2733"
$ws.Range("E37").Value = "  +3.27%  "
$ws.Range("D38").Value = "'1.933"
$ws.Range("E38").Value = "  +1.13%  "
$ws.Range("D39").Value = "'14.43"
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("D40").Value = "'0.09087"
$ws.Range("E40").Value = "  +3.32%  "
$ws.Range("D41").Value = "'0.02732"
$ws.Range("E41").Value = "  +7.91%  "
$ws.Range("D42").Value = "'1.485"
$ws.Range("E42").Value = "  +2.10%  "
$ws.Range("D43").Value = "'0.7658"
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("D44").Value = "'0.7164"
$ws.Range("E44").Value = "  +1.66%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'2.565"
$ws.Range("E45").Value = "  +6.26%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'15.49"
$ws.Range("E46").Value = "  +3.16%  "
$ws.Range("D47").Value = "'4.220"
$ws.Range("E47").Value = "  +2.87%  "
$ws.Range("D48").Value = "'1.001"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").Value = "'140.97"
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("D50").Value = "'1.316"
$ws.Range("E50").Value = "  +8.18%  "
$ws.Range("D51").Value = "'0.07982"
$ws.Range("E51").Value = "  +2.36%  "
